$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9922.429
$ws.Range("I32").Value = 13122.8
$ws.Range("J32").Value = 1921.5
$ws.Range("K32").Value = 13122.8
$ws.Range("L32").Value = 1921.5
$ws.Range("M32").Value = -12796.8
$ws.Range("N32").Value = -2573.5
$ws.Range("H98").Value = 1715.7142
$ws.Range("I98").Value = 1678.7778
$ws.Range("K98").Value = 1678.7778
$ws.Range("M98").Value = -180.7778000000001
$ws.Range("H122").Value = 1715.7142
$ws.Range("I122").Value = 1678.7778
$ws.Range("K122").Value = 5036.3334
$ws.Range("M122").Value = -2586.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29405.727
$ws.Range("I32").Value = 29821.324
$ws.Range("K32").Value = 29821.324
$ws.Range("M32").Value = -29534.324
$ws.Range("H43").Value = 33429.25
$ws.Range("I43").Value = 20342
$ws.Range("J43").Value = 37791.668
$ws.Range("K43").Value = 20342
$ws.Range("L43").Value = 37791.668
$ws.Range("M43").Value = -20029
$ws.Range("N43").Value = -38417.668
$ws.Range("H45").Value = 2554.7856
$ws.Range("I45").Value = 2146.7
$ws.Range("J45").Value = 3575
$ws.Range("K45").Value = 2146.7
$ws.Range("L45").Value = 3575
$ws.Range("M45").Value = -1769.7
$ws.Range("N45").Value = -4329
$ws.Range("H110").Value = 7576967
$ws.Range("I110").Value = 9259985
$ws.Range("J110").Value = 3384.1667
$ws.Range("K110").Value = 9259985
$ws.Range("L110").Value = 3384.1667
$ws.Range("M110").Value = -9257940
$ws.Range("N110").Value = -7474.1667
$ws.Range("H122").Value = 7749.5
$ws.Range("I122").Value = 7749.5
$ws.Range("K122").Value = 23248.5
$ws.Range("M122").Value = -20798.5
$ws.Range("H128").Value = 210197.6
$ws.Range("J128").Value = 210197.6
$ws.Range("L128").Value = 210197.6
$ws.Range("N128").Value = -220157.6
$ws.Range("H132").Value = 4546.698
$ws.Range("I132").Value = 3218.1025
$ws.Range("K132").Value = 9654.307499999999
$ws.Range("M132").Value = -7124.307499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1020.85297
$ws.Range("I94").Value = 876.8095
$ws.Range("J94").Value = 1253.5385
$ws.Range("K94").Value = 876.8095
$ws.Range("L94").Value = 1253.5385
$ws.Range("M94").Value = -425.8095
$ws.Range("N94").Value = -2155.5385
$ws.Range("H105").Value = 45462812
$ws.Range("I105").Value = 90916180
$ws.Range("K105").Value = 90916180
$ws.Range("M105").Value = -90914433

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6291.75
$ws.Range("I58").Value = 6760
$ws.Range("K58").Value = 6760
$ws.Range("M58").Value = -6557
$ws.Range("H134").Value = 35775
$ws.Range("I134").Value = 57000
$ws.Range("K134").Value = 171000
$ws.Range("M134").Value = -168465
$ws.Range("H136").Value = 6291.75
$ws.Range("I136").Value = 6760
$ws.Range("K136").Value = 20280
$ws.Range("M136").Value = -17730

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1315.3334
$ws.Range("I14").Value = 1315.3334
$ws.Range("K14").Value = 3946.0002
$ws.Range("M14").Value = -3773.0002
$ws.Range("H113").Value = 2019.8182
$ws.Range("J113").Value = 2331.6667
$ws.Range("L113").Value = 6995.000100000001
$ws.Range("N113").Value = -11335.0001
$ws.Range("H122").Value = 655.7692
$ws.Range("I122").Value = 399.33334
$ws.Range("J122").Value = 875.5714
$ws.Range("K122").Value = 3594.00006
$ws.Range("L122").Value = 7880.1426
$ws.Range("M122").Value = -1144.00006
$ws.Range("N122").Value = -12780.1426
$ws.Range("H140").Value = 1150
$ws.Range("I140").Value = 795
$ws.Range("J140").Value = 1292
$ws.Range("K140").Value = 2385
$ws.Range("L140").Value = 3876
$ws.Range("M140").Value = 2795
$ws.Range("N140").Value = -14236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5806.4
$ws.Range("I80").Value = 4749.5
$ws.Range("J80").Value = 6511
$ws.Range("K80").Value = 4749.5
$ws.Range("L80").Value = 6511
$ws.Range("M80").Value = -3751.5
$ws.Range("N80").Value = -8507
$ws.Range("H83").Value = 5806.4
$ws.Range("I83").Value = 4749.5
$ws.Range("J83").Value = 6511
$ws.Range("K83").Value = 23747.5
$ws.Range("L83").Value = 32555
$ws.Range("M83").Value = -18755.5
$ws.Range("N83").Value = -42539
$ws.Range("H97").Value = 1068
$ws.Range("I97").Value = 1008.1539
$ws.Range("K97").Value = 1008.1539
$ws.Range("M97").Value = -512.1539
$ws.Range("H113").Value = 2279.1667
$ws.Range("I113").Value = 2193.75
$ws.Range("K113").Value = 2193.75
$ws.Range("M113").Value = -23.75
$ws.Range("H122").Value = 1198.2
$ws.Range("I122").Value = 1166
$ws.Range("J122").Value = 1212
$ws.Range("K122").Value = 3498
$ws.Range("L122").Value = 3636
$ws.Range("M122").Value = -1048
$ws.Range("N122").Value = -8536

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3988.6365
$ws.Range("I22").Value = 2834.25
$ws.Range("J22").Value = 4648.2856
$ws.Range("K22").Value = 2834.25
$ws.Range("L22").Value = 4648.2856
$ws.Range("M22").Value = -2539.25
$ws.Range("N22").Value = -5238.2856
$ws.Range("H27").Value = 3988.6365
$ws.Range("I27").Value = 2834.25
$ws.Range("J27").Value = 4648.2856
$ws.Range("K27").Value = 2834.25
$ws.Range("L27").Value = 4648.2856
$ws.Range("M27").Value = -2727.25
$ws.Range("N27").Value = -4862.2856
$ws.Range("H38").Value = 15995.25
$ws.Range("I38").Value = 15995.5
$ws.Range("K38").Value = 15995.5
$ws.Range("M38").Value = -15585.5
$ws.Range("H40").Value = 9959.799999999999
$ws.Range("I40").Value = 8585.137000000001
$ws.Range("K40").Value = 8585.137000000001
$ws.Range("M40").Value = -8449.137000000001
$ws.Range("H46").Value = 6088.9116
$ws.Range("J46").Value = 6525.0645
$ws.Range("L46").Value = 6525.0645
$ws.Range("N46").Value = -6901.0645
$ws.Range("H82").Value = 1167.5883
$ws.Range("I82").Value = 1191.5834
$ws.Range("J82").Value = 1110
$ws.Range("K82").Value = 1191.5834
$ws.Range("L82").Value = 1110
$ws.Range("M82").Value = -830.5834
$ws.Range("N82").Value = -1832
$ws.Range("H85").Value = 1167.5883
$ws.Range("I85").Value = 1191.5834
$ws.Range("J85").Value = 1110
$ws.Range("K85").Value = 1191.5834
$ws.Range("L85").Value = 1110
$ws.Range("M85").Value = 56.41660000000002
$ws.Range("N85").Value = -3606
$ws.Range("H122").Value = 1835.5
$ws.Range("I122").Value = 1835.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5506.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3056.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 7470.7856
$ws.Range("I132").Value = 5141.875
$ws.Range("K132").Value = 15425.625
$ws.Range("M132").Value = -12895.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 1312.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1312.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1312.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1596.5
$ws.Range("H96").Value = 1299
$ws.Range("I96").Value = 1032.6666
$ws.Range("K96").Value = 1032.6666
$ws.Range("M96").Value = 340.3334
$ws.Range("H122").Value = 2945.182
$ws.Range("I122").Value = 2904.7
$ws.Range("K122").Value = 8714.099999999999
$ws.Range("M122").Value = -6264.099999999999
$ws.Range("H126").Value = 32768.344
$ws.Range("I126").Value = 40915.48
$ws.Range("K126").Value = 122746.44
$ws.Range("M126").Value = -120276.44
$ws.Range("H132").Value = 9446.286
$ws.Range("I132").Value = 7910.8887
$ws.Range("J132").Value = 12210
$ws.Range("K132").Value = 23732.6661
$ws.Range("M132").Value = -21202.6661

